$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PEW")

# Insert 3 new columns before column D (old D "Terms Typically Offered" shifts to G)
$ws.Range("D1:F1").EntireColumn.Insert()

# New header row values
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"

# Fill NA for the new columns in data rows 2-11
$ws.Range("D2:F11").Value = "NA"
